$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item(3)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# --- Paragraph 1: "Abstraction:" -> "Polymorphism:" ---
# Only the word "Abstraction" is retyped as "Polymorphism"; the trailing colon's
# run is left untouched, which is exactly how PowerPoint splits the run when you
# select part of a run and type new text over it.
$para1 = $tr.Paragraphs(1, 1)
$para1.Characters(1, 11).Text = "Polymorphism"

# --- Paragraph 2: reword the bullet under it ---
# "Abstract method startSort(Element[]) is defined in GeneralSort class and
#  implemented in bubble sort, heap sort, and shell sort class"
# becomes
# "Method startSort(Element[]) is defined in GeneralSort class and is overried
#  in bubble sort, heap sort, and shell sort class"
# Addressing is done relative to the paragraph (via a freshly-fetched
# Paragraphs(2,1) range) so it stays correct no matter what happened earlier in
# the text frame. Edits are applied right-to-left within the paragraph so that
# earlier offsets remain valid.
$para2 = $tr.Paragraphs(2, 1)

# Force a run boundary right before "in bubble sort, ..." (no text change).
$r7 = $para2.Characters(86, 47)
$r7.Text = $r7.Text

# Force a run boundary at the single space before "in bubble..." (no text change).
$r6 = $para2.Characters(85, 1)
$r6.Text = $r6.Text

# "lemented" -> "overried"
$para2.Characters(77, 8).Text = "overried"

# "imp" -> "is "
$para2.Characters(74, 3).Text = "is "

# Force a run boundary between "ethod " and "startSort(...)" (no text change).
$r3 = $para2.Characters(17, 57)
$r3.Text = $r3.Text

# "Abstract m" -> "M"
$para2.Characters(1, 10).Text = "M"
